$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")

# Row 11 is a duplicate "Contact / No display for ContactDetail" row.
# Delete it outright so the rows below shift up by one (15 -> 14 total rows),
# matching the new layout where "Contact" is replaced by "Publisher" +
# "Jurisdiction" rows (one row net removed).
$ws1.Rows("11:11").Delete()

# Version bump
$ws1.Range("B3").Value = "6.0.0"

# Date bump
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$ws1.Range("B9").Value = "Alvearie Team"

# The remaining "Contact" row (row 10) becomes "Jurisdiction"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"
